$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2302.1428
$ws.Range("I40").Value = 2098.25
$ws.Range("J40").Value = 2574
$ws.Range("K40").Value = 2098.25
$ws.Range("L40").Value = 2574
$ws.Range("M40").Value = -1923.25
$ws.Range("N40").Value = -2924
$ws.Range("H106").Value = 2709.375
$ws.Range("I106").Value = 3196.8333
$ws.Range("J106").Value = 1247
$ws.Range("K106").Value = 3196.8333
$ws.Range("L106").Value = 1247
$ws.Range("M106").Value = -2565.8333
$ws.Range("N106").Value = -2509
$ws.Range("H113").Value = 24033.445
$ws.Range("J113").Value = 2449
$ws.Range("L113").Value = 2449
$ws.Range("N113").Value = -8957
$ws.Range("H116").Value = 3996.1428
$ws.Range("I116").Value = 4000
$ws.Range("J116").Value = 3995.5
$ws.Range("K116").Value = 4000
$ws.Range("L116").Value = 3995.5
$ws.Range("M116").Value = -558
$ws.Range("N116").Value = -10879.5
$ws.Range("H125").Value = 581.6667
$ws.Range("I125").Value = 607.375
$ws.Range("K125").Value = 5466.375
$ws.Range("M125").Value = -3006.375
$ws.Range("H132").Value = 1095.8148
$ws.Range("I132").Value = 984.5238000000001
$ws.Range("K132").Value = 2953.5714
$ws.Range("M132").Value = -423.5714000000003
$ws.Range("H137").Value = 1983.1666
$ws.Range("I137").Value = 1249.75
$ws.Range("K137").Value = 3749.25
$ws.Range("M137").Value = -1199.25

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2777853
$ws.Range("I2").Value = 2777853
$ws.Range("K2").Value = 2777853
$ws.Range("M2").Value = -2777740
$ws.Range("H32").Value = 6119.6665
$ws.Range("I32").Value = 5082.9116
$ws.Range("J32").Value = 23744.5
$ws.Range("K32").Value = 5082.9116
$ws.Range("L32").Value = 23744.5
$ws.Range("M32").Value = -4795.9116
$ws.Range("N32").Value = -24318.5
$ws.Range("H61").Value = 5361.2964
$ws.Range("I61").Value = 6643.3125
$ws.Range("K61").Value = 6643.3125
$ws.Range("M61").Value = -6431.3125
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
$ws.Range("H102").Value = 1352.8889
$ws.Range("I102").Value = 1246.5
$ws.Range("J102").Value = 1438
$ws.Range("K102").Value = 1246.5
$ws.Range("L102").Value = 1438
$ws.Range("M102").Value = 375.5
$ws.Range("N102").Value = -4682
$ws.Range("H109").Value = 58659.25
$ws.Range("J109").Value = 58659.25
$ws.Range("L109").Value = 58659.25
$ws.Range("N109").Value = -61433.25
$ws.Range("H116").Value = 2777853
$ws.Range("I116").Value = 2777853
$ws.Range("K116").Value = 2777853
$ws.Range("M116").Value = -2775559
$ws.Range("H122").Value = 1413.7222
$ws.Range("I122").Value = 1072.2727
$ws.Range("J122").Value = 1950.2858
$ws.Range("K122").Value = 3216.8181
$ws.Range("L122").Value = 5850.857400000001
$ws.Range("M122").Value = -766.8181
$ws.Range("N122").Value = -10750.8574
$ws.Range("H132").Value = 1627.6333
$ws.Range("I132").Value = 1529.7142
$ws.Range("J132").Value = 2998.5
$ws.Range("K132").Value = 4589.142599999999
$ws.Range("L132").Value = 8995.5
$ws.Range("M132").Value = -2059.142599999999
$ws.Range("N132").Value = -14055.5
$ws.Range("H136").Value = 5361.2964
$ws.Range("I136").Value = 6643.3125
$ws.Range("K136").Value = 19929.9375
$ws.Range("M136").Value = -17379.9375

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2777853
$ws.Range("I3").Value = 2777853
$ws.Range("K3").Value = 2777853
$ws.Range("M3").Value = -2777739
$ws.Range("H135").Value = 52999.5
$ws.Range("J135").Value = 52999.5
$ws.Range("L135").Value = 52999.5
$ws.Range("N135").Value = -63139.5
$ws.Range("H140").Value = 32145
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 32145
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 32145
$ws.Range("M140").ClearContents()
$ws.Range("N140").Value = -42505

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1501050.9
$ws.Range("I58").Value = 2558697.5
$ws.Range("J58").Value = 2718.3333
$ws.Range("K58").Value = 2558697.5
$ws.Range("L58").Value = 2718.3333
$ws.Range("M58").Value = -2558494.5
$ws.Range("N58").Value = -3124.3333
$ws.Range("H99").Value = 1796.4
$ws.Range("I99").Value = 1995.5
$ws.Range("J99").Value = 1000
$ws.Range("K99").Value = 1995.5
$ws.Range("L99").Value = 1000
$ws.Range("M99").Value = -497.5
$ws.Range("N99").Value = -3996
$ws.Range("H126").Value = 1796.4
$ws.Range("I126").Value = 1995.5
$ws.Range("J126").Value = 1000
$ws.Range("K126").Value = 5986.5
$ws.Range("L126").Value = 3000
$ws.Range("M126").Value = -3516.5
$ws.Range("N126").Value = -7940
$ws.Range("H132").Value = 1828.65
$ws.Range("I132").Value = 1149.7059
$ws.Range("J132").Value = 5676
$ws.Range("K132").Value = 3449.1177
$ws.Range("L132").Value = 17028
$ws.Range("M132").Value = -919.1176999999998
$ws.Range("N132").Value = -22088
$ws.Range("H134").Value = 1542.7241
$ws.Range("I134").Value = 1529.4642
$ws.Range("K134").Value = 4588.392599999999
$ws.Range("M134").Value = -2053.392599999999
$ws.Range("H136").Value = 1501050.9
$ws.Range("I136").Value = 2558697.5
$ws.Range("J136").Value = 2718.3333
$ws.Range("K136").Value = 7676092.5
$ws.Range("L136").Value = 8154.999899999999
$ws.Range("M136").Value = -7673542.5
$ws.Range("N136").Value = -13254.9999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 476
$ws.Range("I5").Value = 443.3846
$ws.Range("K5").Value = 1330.1538
$ws.Range("M5").Value = -1218.1538
$ws.Range("H114").Value = 2376
$ws.Range("J114").Value = 3569.8333
$ws.Range("L114").Value = 10709.4999
$ws.Range("N114").Value = -17217.4999
$ws.Range("H129").Value = 60852.582
$ws.Range("I129").Value = 724.25
$ws.Range("J129").Value = 90916.75
$ws.Range("K129").Value = 2172.75
$ws.Range("L129").Value = 272750.25
$ws.Range("M129").Value = 2827.25
$ws.Range("N129").Value = -282750.25
$ws.Range("H131").Value = 17048.488
$ws.Range("J131").Value = 18724.486
$ws.Range("L131").Value = 56173.458
$ws.Range("N131").Value = -66253.458
$ws.Range("H132").Value = 2403.1667
$ws.Range("I132").Value = 1220
$ws.Range("J132").Value = 2639.8
$ws.Range("K132").Value = 10980
$ws.Range("L132").Value = 23758.2
$ws.Range("M132").Value = -8450
$ws.Range("N132").Value = -28818.2
$ws.Range("H135").Value = 476
$ws.Range("I135").Value = 443.3846
$ws.Range("K135").Value = 3990.4614
$ws.Range("M135").Value = -1455.4614
$ws.Range("H140").Value = 1849.6666
$ws.Range("I140").Value = 859.76
$ws.Range("K140").Value = 2579.28
$ws.Range("M140").Value = 2600.72

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 54304.5
$ws.Range("I19").Value = 200
$ws.Range("K19").Value = 200
$ws.Range("M19").Value = 88
$ws.Range("H80").Value = 2553.875
$ws.Range("J80").Value = 2000
$ws.Range("L80").Value = 2000
$ws.Range("N80").Value = -3996
$ws.Range("H83").Value = 2553.875
$ws.Range("J83").Value = 2000
$ws.Range("L83").Value = 10000
$ws.Range("N83").Value = -19984
$ws.Range("H126").Value = 2926811.8
$ws.Range("I126").Value = 4276418
$ws.Range("K126").Value = 12829254
$ws.Range("M126").Value = -12826784
$ws.Range("H132").Value = 1480362.1
$ws.Range("I132").Value = 1749309.9
$ws.Range("J132").Value = 1149.5
$ws.Range("K132").Value = 5247929.699999999
$ws.Range("L132").Value = 3448.5
$ws.Range("M132").Value = -5245399.699999999
$ws.Range("N132").Value = -8508.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1655.625
$ws.Range("I46").Value = 1014.3333
$ws.Range("J46").Value = 2480.1428
$ws.Range("K46").Value = 1014.3333
$ws.Range("L46").Value = 2480.1428
$ws.Range("M46").Value = -826.3333
$ws.Range("N46").Value = -2856.1428
$ws.Range("H122").Value = 14390.1
$ws.Range("I122").Value = 13737.625
$ws.Range("K122").Value = 41212.875
$ws.Range("M122").Value = -38762.875

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()
$ws.Range("H126").Value = 4807.75
$ws.Range("I126").Value = 6757.4443
$ws.Range("J126").Value = 1298.3
$ws.Range("K126").Value = 20272.3329
$ws.Range("L126").Value = 3894.9
$ws.Range("M126").Value = -17802.3329
$ws.Range("N126").Value = -8834.9
$ws.Range("H132").Value = 1928.4375
$ws.Range("I132").Value = 1326.5
$ws.Range("J132").Value = 3734.25
$ws.Range("K132").Value = 3979.5
$ws.Range("L132").Value = 11202.75
$ws.Range("M132").Value = -1449.5
$ws.Range("N132").Value = -16262.75
